# API Automation Framework Part-V Added Logging for CowinAPI TC's
#
# On the TESTDATA sheet, the test row for the Kanpur Nagar "Yes" execute
# case (row 6) gets a new pin code and a new expected output value:
#   C6: 226003 -> 226012
#   D6: "HWC NEWAL BANGARMAU" -> "Apollo"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTDATA")

$ws.Range("C6").Value = 226012
$ws.Range("D6").Value = "Apollo"
